# Correcciones en reglas del documento stock actual
# Applies the "Diferencia Stock" / "Stock Real" / "Pedido Corregido Stock" /
# "Ventas Objetivo" / "Beneficio Objetivo" / "Pedido Final" corrections for
# Semana_6, plus the resulting summary-metric updates and the re-hiding of
# the rows whose corrected "Pedido Corregido Stock" (column Q) dropped to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0) Rows 13 and 23 are already hidden in the source and stay hidden - only a
#    cell inside them changes. Temporarily unhide before writing so the
#    engine doesn't stamp a spurious explicit row height on a write into a
#    zero-height hidden row; the real hidden rows get re-hidden in step 3.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Hidden = $false
$ws.Rows.Item(23).Hidden = $false

# ---------------------------------------------------------------------------
# 1) Plain numeric cell value corrections (column L/M/N/P/Q/U per row).
#    Kept as one big table so every (cell -> new value) pair is explicit and
#    easy to audit against the source diff.
# ---------------------------------------------------------------------------
$cellValues = @{
    "L13" = 1

    "M16" = 0
    "N16" = 0
    "P16" = 16
    "Q16" = 0
    "U16" = 0

    "M17" = 0
    "N17" = 0
    "P17" = 9
    "Q17" = 0
    "U17" = 0

    "M18" = 0
    "N18" = 0
    "P18" = 24
    "Q18" = 0
    "U18" = 0

    "L19" = 1
    "M19" = 0
    "N19" = 0
    "P19" = 8
    "Q19" = 0
    "U19" = 0

    "M20" = 0
    "N20" = 0
    "P20" = 21
    "Q20" = 0
    "U20" = 0

    "L23" = 1

    "M27" = 75.90000000000001
    "N27" = 45.54
    "P27" = 2
    "Q27" = 4
    "U27" = 4

    "L29" = 1
    "M29" = 0
    "N29" = 0
    "P29" = 10
    "Q29" = 0
    "U29" = 0

    "L30" = 1

    "M32" = 0
    "N32" = 0
    "P32" = 10
    "Q32" = 0
    "U32" = 0

    "L33" = 1

    "M34" = 43.25
    "N34" = 25.95
    "P34" = 1
    "Q34" = 2
    "U34" = 2

    "L35" = 1
    "M35" = 0
    "N35" = 0
    "P35" = 3
    "Q35" = 0
    "U35" = 0

    "L37" = 1
    "M37" = 21.65
    "N37" = 12.99
    "P37" = 2
    "Q37" = 1
    "U37" = 1

    "L39" = 1

    "L40" = 1

    "M41" = 19.45
    "N41" = 11.67
    "P41" = 2
    "Q41" = 1
    "U41" = 1

    "M45" = 56.75
    "N45" = 34.05
    "P45" = 1
    "Q45" = 1
    "U45" = 1

    "M46" = 0
    "N46" = 0
    "P46" = 10
    "Q46" = 0
    "U46" = 0

    "C50" = 62
    "C61" = 10
}

foreach ($cellRef in $cellValues.Keys) {
    $ws.Range($cellRef).Value = $cellValues[$cellRef]
}

# ---------------------------------------------------------------------------
# 2) Total_Importe summary text (keeps the € suffix as literal text, not as
#    an auto-converted currency number, so the cell keeps its original
#    General-format style instead of Excel inferring a new currency format).
# ---------------------------------------------------------------------------
$ws.Range("C52").Formula = '="1244.16€"'
$ws.Range("C52").Copy()
$ws.Range("C52").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Rows whose corrected "Pedido Corregido Stock" (Q) fell to 0 are now
#    hidden (they no longer need ordering this week), and rows 13/23 return
#    to their original hidden state.
# ---------------------------------------------------------------------------
$rowsToHide = @(13, 16, 17, 18, 19, 20, 23, 29, 32, 35, 46)
foreach ($rowNum in $rowsToHide) {
    $ws.Rows.Item($rowNum).Hidden = $true
}
